# Updated cryptos list on Sun Apr 23 20:52:05 UTC 2023 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) figures for
# each coin row to match the latest upstream data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row number -> @{ D = "new price text"; E = "new volume text" }
# (only the columns that actually changed are present per row)
$updates = @{
    2 = @{ D = "27.921.76"; E = "  +0.61%  " }
    3 = @{ D = "1.878.48"; E = "  -0.01%  " }
    4 = @{ E = "  +1.45%  " }
    5 = @{ D = "334.77"; E = "  +0.84%  " }
    6 = @{ E = "  +1.30%  " }
    7 = @{ D = "0.4694"; E = "  -0.76%  " }
    8 = @{ D = "0.3909"; E = "  -1.33%  " }
    9 = @{ D = "46.87"; E = "  -2.04%  " }
    10 = @{ D = "0.07946"; E = "  -0.96%  " }
    11 = @{ E = "  -1.67%  " }
    12 = @{ D = "21.61"; E = "  -1.10%  " }
    13 = @{ D = "1.890.78"; E = "  +0.25%  " }
    14 = @{ D = "5.948"; E = "  -0.28%  " }
    15 = @{ D = "7.097"; E = "  -0.89%  " }
    16 = @{ E = "  +1.59%  " }
    17 = @{ D = "0.06790"; E = "  +2.50%  " }
    18 = @{ D = "87.57"; E = "  +0.42%  " }
    19 = @{ E = "  -0.11%  " }
    20 = @{ D = "17.04"; E = "  -1.64%  " }
    21 = @{ E = "  +1.34%  " }
    22 = @{ D = "27.912.12"; E = "  +0.52%  " }
    23 = @{ D = "5.468"; E = "  -0.78%  " }
    24 = @{ D = "10.97"; E = "  -0.65%  " }
    25 = @{ D = "2.362"; E = "  +2.75%  " }
    26 = @{ D = "2.100.08"; E = "  -0.42%  " }
    27 = @{ D = "159.67"; E = "  +2.08%  " }
    28 = @{ E = "  -1.75%  " }
    29 = @{ D = "2.075"; E = "  -1.12%  " }
    30 = @{ D = "5.451"; E = "  -2.54%  " }
    31 = @{ D = "120.92"; E = "  -1.38%  " }
    32 = @{ E = "  -0.33%  " }
    33 = @{ D = "0.9569"; E = "  -1.23%  " }
    34 = @{ D = "3.658"; E = "  +0.89%  " }
    35 = @{ D = "5.310"; E = "  +0.11%  " }
    36 = @{ D = "1.352"; E = "  -7.15%  " }
    37 = @{ D = "0.06105"; E = "  -0.18%  " }
    38 = @{ D = "0.02242" }
    39 = @{ D = "1.204"; E = "  -1.96%  " }
    40 = @{ E = "  +1.32%  " }
    41 = @{ D = "8.121"; E = "  -1.02%  " }
    42 = @{ D = "0.5891"; E = "  -1.64%  " }
    43 = @{ D = "0.1892"; E = "  -1.05%  " }
    44 = @{ D = "10.21"; E = "  -0.39%  " }
    45 = @{ D = "1.271"; E = "  +1.84%  " }
    46 = @{ D = "0.5645"; E = "  -1.26%  " }
    47 = @{ D = "12.15"; E = "  -1.14%  " }
    48 = @{ D = "3.395"; E = "  -0.25%  " }
    49 = @{ D = "1.918"; E = "  -0.81%  " }
    50 = @{ D = "0.06861"; E = "  +0.59%  " }
    51 = @{ D = "113.53"; E = "  +1.20%  " }
}

foreach ($r in $updates.Keys) {
    $vals = $updates[$r]
    if ($vals.ContainsKey("D")) {
        # Keep the cell text even though the new price string often looks
        # like a plain number ("21.61"); NumberFormat "@" forces text entry
        # and the Style reset back to "Normal" drops the temporary format
        # so the saved cell ends up with no explicit style, same as source.
        $ws.Range("D$r").NumberFormat = "@"
        $ws.Range("D$r").Value = $vals["D"]
        $ws.Range("D$r").Style = "Normal"
    }
    if ($vals.ContainsKey("E")) {
        $ws.Range("E$r").NumberFormat = "@"
        $ws.Range("E$r").Value = $vals["E"]
        $ws.Range("E$r").Style = "Normal"
    }
}

Write-Output "Applied cryptos price/volume refresh"
